# 27-12-2017 : EmpMast Import removed MedCheckFlg and PoliceVerificationFlg
#
# The upload-template sheet ("Sheet1") had two columns - "MedChkFlg" (O)
# and "PoliceVeriFlg" (P) - that are no longer required by the import, so
# they are removed entirely (cells, column widths and all cells to their
# right shift left by two columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the MedChkFlg (O) and PoliceVeriFlg (P) columns in one operation
# so everything after them (ShiftCode, EmpCode, ContCode, ... ValidTo)
# shifts two columns to the left.
$ws.Range("O1:P1").EntireColumn.Delete() | Out-Null

# Leave the same kind of cell selected/active as in the saved workbook.
$ws.Range("O15").Select() | Out-Null
